$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The NIK column (C2:C7) all point at the same "EN-4-025" entry, so every one
# of those cells has to be rewritten to "EN-4-095" together - otherwise only
# C2 would change and the old string would linger for C3:C7.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 3).Value = "EN-4-095"
}

# Employee name for the first row changes too.
$ws.Range("D2").Value = "Rahmat Hidayat"

# Running index in column A is renumbered from 1-6 to 7-12.
$newIndex = @(7, 8, 9, 10, 11, 12)
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $newIndex[$r - 2]
}

# Move the active selection to D8, matching the saved cursor position.
$ws.Range("D8").Select()

$wb.Save()
